$p = $ppt.ActivePresentation

# --- Slide 1 (title slide): set title + subtitle text ---
$s1 = $p.Slides.Item(1)

$titleTr = $s1.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "EudAssistent"
$titleTr.InsertAfter(" Demo")

$s1.Shapes.Item(2).TextFrame.TextRange.Text = "Vlad Fernoaga"

# --- Slide 2 (existing quiz slide): update title + hyperlink text/tooltip ---
$s2 = $p.Slides.Item(2)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Check your knowledge in History"

$s2ContentTr = $s2.Shapes.Item(2).TextFrame.TextRange
$s2ContentTr.Text = "Trigger Question"
$s2ActionSetting = $s2ContentTr.ActionSettings.Item(1)
$s2ActionSetting.Hyperlink.ScreenTip = ""

# --- New slides: duplicate slide 2 (same layout/shapes/hyperlink) four times, ---
# --- then re-point title/content text for each subject.                      ---
# Duplicating slide 2 always re-inserts right after slide 2 (position 3),
# so doing this four times in a row yields slides in id order 261,260,259,258
# at positions 3..6; moving the first (id 261 / "engineering") to the end
# reproduces the authored order 260,259,258,261.

$dHistory     = $p.Slides.Item(2).Duplicate()
$dDemography  = $p.Slides.Item(2).Duplicate()
$dGeography   = $p.Slides.Item(2).Duplicate()
$dEngineering = $p.Slides.Item(2).Duplicate()

$dEngineering.MoveTo(6)

# History slide
$dHistory.Shapes.Item(1).TextFrame.TextRange.Text = "Check your knowledge in History "
$trH = $dHistory.Shapes.Item(2).TextFrame.TextRange
$trH.Text = "Trigger Question"
$trH.ActionSettings.Item(1).Hyperlink.ScreenTip = ""

# Demography slide
$dDemography.Shapes.Item(1).TextFrame.TextRange.Text = "Check your knowledge in Demography"
$trD = $dDemography.Shapes.Item(2).TextFrame.TextRange
$trD.Text = "Trigger Question"
$trD.ActionSettings.Item(1).Hyperlink.ScreenTip = ""

# Geography slide
$dGeography.Shapes.Item(1).TextFrame.TextRange.Text = "Check your knowledge in Geography"
$trG = $dGeography.Shapes.Item(2).TextFrame.TextRange
$trG.Text = "Trigger Question"
$trG.ActionSettings.Item(1).Hyperlink.ScreenTip = ""

# Engineering slide
$dEngineering.Shapes.Item(1).TextFrame.TextRange.Text = "Check your knowledge in engineering"
$trE = $dEngineering.Shapes.Item(2).TextFrame.TextRange
$trE.Text = "Trigger Question"
$trE.ActionSettings.Item(1).Hyperlink.ScreenTip = ""
